$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 100
$ws.Range("F5").Value = 5890
$ws.Range("F6").Value = 486
$ws.Range("F7").Value = 1049
$ws.Range("F8").Value = 11
$ws.Range("F9").Value = 3442
$ws.Range("F10").Value = 6668
$ws.Range("F11").Value = 209
$ws.Range("F12").Value = 1318
$ws.Range("F13").Value = 764
$ws.Range("F16").Value = 24
$ws.Range("F17").Value = 1117
$ws.Range("F21").Value = 175
$ws.Range("F23").Value = 978
$ws.Range("F26").Value = 16
$ws.Range("F29").Value = 1154
$ws.Range("F31").Value = 54
$ws.Range("F35").Value = 307
$ws.Range("F36").Value = 16
$ws.Range("F37").Value = 47
$ws.Range("F38").Value = 292
$ws.Range("F39").Value = 1168
$ws.Range("F40").Value = 53
$ws.Range("F41").Value = 98

$ws = $wb.Worksheets.Item(2)
$ws.Range("F10").Value = 30
$ws.Range("F26").Value = 609
$ws.Range("F30").Value = 684
$ws.Range("F32").Value = 570
$ws.Range("F38").Value = 124
$ws.Range("F40").Value = 53

$ws = $wb.Worksheets.Item(3)
$ws.Range("F8").Value = 1058

$ws = $wb.Worksheets.Item(4)
$ws.Range("F7").Value = 100
$ws.Range("F14").Value = 5890
$ws.Range("F15").Value = 486
$ws.Range("F16").Value = 1049
$ws.Range("F17").Value = 3442
$ws.Range("F18").Value = 30
$ws.Range("F19").Value = 6668
$ws.Range("F20").Value = 209
$ws.Range("F21").Value = 1318
$ws.Range("F24").Value = 764
$ws.Range("F26").Value = 1058
$ws.Range("F28").Value = 24
$ws.Range("F29").Value = 1117
$ws.Range("F31").Value = 175
$ws.Range("F32").Value = 978
$ws.Range("F33").Value = 609
$ws.Range("F35").Value = 16
$ws.Range("F37").Value = 1154
$ws.Range("F39").Value = 54
$ws.Range("F43").Value = 570
$ws.Range("F44").Value = 307
$ws.Range("F46").Value = 292
$ws.Range("F48").Value = 124
$ws.Range("F50").Value = 98
$ws.Range("F51").Value = 53
